# Update the "取得日時" (acquired datetime) timestamps in rows 2-6 of the
# first worksheet ("ランサーズ") from 2025-12-28 12:36:26 to 2025-12-28 12:48:23.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$newTimestamp = "2025-12-28 12:48:23"

for ($row = 2; $row -le 6; $row++) {
    $ws.Cells.Item($row, 1).Value = $newTimestamp
}
